# Daily attendance processing - 2026-01-13 10:06:16
# Reorders the "Recorded By" (column G) contributor list on each data row so
# that any "system"/"System" entries come first, followed by the remaining
# entries (e.g. user emails) in their original relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ([string]::IsNullOrEmpty($value)) {
        continue
    }

    $parts = $value -split ", "
    if ($parts.Count -le 1) {
        continue
    }

    $systemParts = @()
    $otherParts = @()
    foreach ($part in $parts) {
        if ($part.ToLower().Contains("system")) {
            $systemParts += $part
        } else {
            $otherParts += $part
        }
    }

    $newParts = $systemParts + $otherParts
    $newValue = $newParts -join ", "

    if ($newValue -ne $value) {
        $cell.Value = $newValue
    }
}
